# TC05_Bento_MultiFilter_Diagnosis-Recurrence-TumorSize-Chemo-ERStatus.xlsx
# "startup" sheet: the four Cypher scripts (one per tab: CasesTab/B2,
# SamplesTab/B3, FilesTab/B4, and the shared StatQuery in column C) all
# filter on d.er_status. Flip that filter from "Negative" to "Positive"
# in every one of them (updated automation scripts for the Bento perf
# pass), and grow the row heights slightly to fit the rewritten text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

function Set-ErStatusPositive($rangeAddress) {
    $cell = $ws.Range($rangeAddress)
    $current = $cell.Value()
    $updated = $current.Replace('d.er_status In ["Negative"]', 'd.er_status In ["Positive"]')
    $cell.Value = $updated
}

# CasesTab query (row 2)
Set-ErStatusPositive "B2"
# SamplesTab query (row 3)
Set-ErStatusPositive "B3"
# FilesTab query (row 4)
Set-ErStatusPositive "B4"

# Shared StatQuery text, duplicated across C2:C4
Set-ErStatusPositive "C2"
Set-ErStatusPositive "C3"
Set-ErStatusPositive "C4"

# Row heights grew a touch to fit the (same-length) rewritten text
$ws.Rows.Item(2).RowHeight = 375
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 330

# Saved view now shows row 4 (FilesTab) scrolled into place with C4 active
$ws.Range("C4").Select() | Out-Null
